$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.861.48"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "'2.922.61"
$ws.Range("E3").Value = "  +3.47%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'352.32"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").Value = "'111.88"
$ws.Range("E6").Value = "  -0.92%  "

$ws.Range("D7").Value = "'0.561"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").Value = "'39.30"
$ws.Range("E10").Value = "  -1.98%  "

$ws.Range("D11").Value = "'0.0878"
$ws.Range("E11").Value = "  +3.30%  "

$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("D13").Value = "'20.09"
$ws.Range("E13").Value = "  +0.57%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'7.76"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "'3.384.56"
$ws.Range("E15").Value = "  +3.43%  "

$ws.Range("D16").Value = "'2.925.88"
$ws.Range("E16").Value = "  +3.64%  "

$ws.Range("D17").Value = "'0.982"
$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").Value = "'51.917.48"
$ws.Range("E18").Value = "  +0.27%  "

$ws.Range("B19").Value = "ImmutableX"
$ws.Range("C19").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D19").Value = "'3.31"
$ws.Range("E19").Value = "  -3.91%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'7.60"
$ws.Range("E20").Value = "  +0.09%  "

$ws.Range("D21").Value = "'14.22"
$ws.Range("E21").Value = "  +6.63%  "

$ws.Range("D22").Value = "'0.0₃0979"
$ws.Range("E22").Value = "  +0.63%  "

$ws.Range("D23").Value = "'71.20"
$ws.Range("E23").Value = "  +0.93%  "

$ws.Range("D24").Value = "'268.31"
$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("E25").Value = "  +0.73%  "

$ws.Range("E26").Value = "  +11.32%  "

$ws.Range("D27").Value = "'26.93"
$ws.Range("E27").Value = "  +2.59%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.14%  "

$ws.Range("E29").Value = "  +17.14%  "

$ws.Range("E30").Value = "  +17.14%  "

$ws.Range("D31").Value = "'10.58"
$ws.Range("E31").Value = "  +0.07%  "

$ws.Range("D32").Value = "'37.19"
$ws.Range("E32").Value = "  -4.96%  "

$ws.Range("E33").Value = "  +0.20%  "

$ws.Range("E34").Value = "  +10.09%  "

$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("D36").Value = "'0.0453"
$ws.Range("E36").Value = "  -1.33%  "

$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("D38").Value = "'3.33"
$ws.Range("E38").Value = "  +3.34%  "

$ws.Range("D39").Value = "'18.65"
$ws.Range("E39").Value = "  -2.13%  "

$ws.Range("D40").Value = "'2.05"
$ws.Range("E40").Value = "  +1.72%  "

$ws.Range("E41").Value = "  +6.36%  "

$ws.Range("E42").Value = "  +1.49%  "

$ws.Range("D43").Value = "'23.39"
$ws.Range("E43").Value = "  +6.20%  "

$ws.Range("E44").Value = "  -1.30%  "

$ws.Range("E45").Value = "  +2.02%  "

$ws.Range("D46").Value = "'3.51"
$ws.Range("E46").Value = "  +0.58%  "

$ws.Range("D47").Value = "'2.170.71"
$ws.Range("E47").Value = "  -0.63%  "

$ws.Range("D48").Value = "'110.99"
$ws.Range("E48").Value = "  -8.79%  "

$ws.Range("E49").Value = "  -0.13%  "

$ws.Range("E50").Value = "  +8.85%  "

$ws.Range("D51").Value = "'0.946"
$ws.Range("E51").Value = "  -4.25%  "
